# Add "Management Info Manager Name" column (AQ) to the
# ApprovedProcessData-Morn sheet, filled down through every existing
# data row plus one trailing blank row, matching the justified /
# wrap-text header style used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerText = "Management Info Manager Name"

# Last row that should receive the new column: one past the current
# last used row (44 data rows -> 45 including the new trailing row).
$lastRow = $ws.UsedRange.Rows.Count + 1

$col = $ws.Range("AQ1:AQ" + $lastRow)
$col.Value = $headerText

# Style to match the existing header formatting (bold, justified,
# wrap text) - mirrors the xf used for column headers elsewhere.
# NOTE: order matters for how the style table gets built up, so
# WrapText is applied before the alignment/bold tweaks.
$col.WrapText = $true
$col.HorizontalAlignment = -4130
$col.VerticalAlignment = -4130
$col.Font.Bold = $true

# Row heights grow to fit the new wrapped header text.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Rows.Item($r).RowHeight = 59.15
}

$ws.Range("AQ2").Select()
